$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold text that looks numeric (e.g. "1.005", "30.014.64").
# Force text entry (quotePrefix) then reset the style back to Normal so the
# saved cell keeps the original (style-less) text-string representation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.995.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.894.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.38%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4962'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.89'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2963'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06643'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.900.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07251'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6779'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.862'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.016.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007985'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.148.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.774'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.14%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.670'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.06%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.259'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '148.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '131.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.961'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.383'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.219'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08749'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.940'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05095'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.122'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7028'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.685'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.790'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.223'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9555'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01664'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.986'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.46%  '
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4235'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.449'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.36%  '
$ws.Range("E47").Value = '  +3.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05753'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3740'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.77%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.200'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.18%  '
